# "run some of the 27 and 31 degree TA"
#
# Adds the CRM-accuracy titration rows that were run for the 27 and 31
# degree TA batches: fills in the date/batch-value that had been missing
# for row 10, and appends rows 11-16 recording further runs (row 13 is
# marked as the end of the sample).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$crm        = 2207.0300000000002   # reference CRM value used throughout
$sampleDate = 43187                 # 3/28/2018
$batchNum   = 169

# --- Row 10 already existed as a CRM-only placeholder row; fill in the
#     date and measured batch value that was recorded for it ---
$ws.Range("A2").Copy()
$ws.Range("A10").PasteSpecial(-4122)   # xlPasteFormats - reuse A2's date style
$ws.Range("A10").Value = $sampleDate
$ws.Range("B10").Value = 2219.4055706962999

# --- Row 11: new sample run ---
$ws.Range("A2").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").Value = $sampleDate
$ws.Range("B11").Value = 2211
$ws.Range("C11").Value = $crm
$ws.Range("D11").Formula = "=100*(B11-C11)/C11"
$ws.Range("E11").Value = $batchNum

# --- Row 12: new sample run ---
$ws.Range("A2").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").Value = $sampleDate
$ws.Range("B12").Value = 2217
$ws.Range("C12").Value = $crm
$ws.Range("E12").Value = $batchNum

# --- Row 13: date recorded but no batch value - marks the end of sample ---
$ws.Range("A2").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = $sampleDate
$ws.Range("C13").Value = $crm
$ws.Range("E13").Value = $batchNum
$ws.Range("F13").Value = "end of sample"

# --- Rows 14-16: trailing rows with just the CRM reference value recorded ---
foreach ($r in 14..16) {
    $ws.Range("C$r").Value = $crm
    $ws.Range("E$r").Value = $batchNum
}

# Fill the "% off" formula down through the newly-added rows 12-16 in one
# shot so that it is stored as a single shared-formula block (D11 keeps its
# own literal formula entered above, same as in the source data).
$ws.Range("D12:D16").Formula = "=100*(B12-C12)/C12"

# --- Selection reflects where the user left off entering data ---
$ws.Range("B13").Select()
